# Insert a new weekly data row for "Femacal de La Calera - Espinaca" before
# the existing row 384. This pushes the previous rows 384-409 down to
# 385-410 (their values are left untouched) and grows the used range from
# A1:R409 to A1:R410.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(384).EntireRow.Insert()

$ws.Cells.Item(384, 1).Value  = 3
$ws.Cells.Item(384, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(384, 3).Value  = "Coquimbo"
$ws.Cells.Item(384, 4).Value  = 44826
$ws.Cells.Item(384, 5).Value  = 5
$ws.Cells.Item(384, 6).Value  = 100112012
$ws.Cells.Item(384, 7).Value  = "Espinaca"
$ws.Cells.Item(384, 8).Value  = "Sin especificar"
$ws.Cells.Item(384, 9).Value  = "Primera"
$ws.Cells.Item(384, 10).Value = 120
$ws.Cells.Item(384, 11).Value = 4000
$ws.Cells.Item(384, 12).Value = 4000
$ws.Cells.Item(384, 13).Value = 4000
$ws.Cells.Item(384, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(384, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(384, 16).Value = 1333
$ws.Cells.Item(384, 17).Value = 3
$ws.Cells.Item(384, 18).Value = "Hortaliza"
